$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.537.89"
$ws.Range("E2").Value = "  -1.06%  "

# Row 3
$ws.Range("D3").Value = "1.592.68"
$ws.Range("E3").Value = "  -1.54%  "

# Row 4
$ws.Range("E4").Value = "  +0.26%  "

# Row 5
$ws.Range("D5").Value = "'207.41"
$ws.Range("E5").Value = "  -1.26%  "

# Row 6
$ws.Range("D6").Value = "'0.499"
$ws.Range("E6").Value = "  -4.06%  "

# Row 7
$ws.Range("E7").Value = "  +0.27%  "

# Row 8
$ws.Range("D8").Value = "'22.22"
$ws.Range("E8").Value = "  -4.33%  "

# Row 9
$ws.Range("E9").Value = "  -1.74%  "

# Row 10
$ws.Range("E10").Value = "  -3.44%  "

# Row 11
$ws.Range("D11").Value = "'0.0870"
$ws.Range("E11").Value = "  -0.80%  "

# Row 12
$ws.Range("D12").Value = "1.821.03"
$ws.Range("E12").Value = "  -1.39%  "

# Row 13
$ws.Range("D13").Value = "1.592.68"
$ws.Range("E13").Value = "  -1.47%  "

# Row 14
$ws.Range("D14").Value = "'3.85"
$ws.Range("E14").Value = "  -4.02%  "

# Row 15
$ws.Range("D15").Value = "'0.537"
$ws.Range("E15").Value = "  -3.85%  "

# Row 16
$ws.Range("D16").Value = "'63.24"
$ws.Range("E16").Value = "  -2.80%  "

# Row 17
$ws.Range("D17").Value = "27.532.07"
$ws.Range("E17").Value = "  -1.01%  "

# Row 18
$ws.Range("D18").Value = "'216.70"
$ws.Range("E18").Value = "  -5.06%  "

# Row 19
$ws.Range("D19").Value = "'7.36"
$ws.Range("E19").Value = "  -3.04%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0691"
$ws.Range("E20").Value = "  -3.83%  "

# Row 21
$ws.Range("E21").Value = "  +0.30%  "

# Row 22
$ws.Range("E22").Value = "  -2.41%  "

# Row 23
$ws.Range("D23").Value = "'9.66"
$ws.Range("E23").Value = "  -3.91%  "

# Row 24
$ws.Range("E24").Value = "  -1.37%  "

# Row 25
$ws.Range("D25").Value = "'155.02"
$ws.Range("E25").Value = "  +0.45%  "

# Row 26
$ws.Range("E26").Value = "  +0.29%  "

# Row 27
$ws.Range("D27").Value = "'6.69"
$ws.Range("E27").Value = "  -2.54%  "

# Row 28
$ws.Range("D28").Value = "'14.98"
$ws.Range("E28").Value = "  -2.93%  "

# Row 29
$ws.Range("E29").Value = "  -4.71%  "

# Row 30
$ws.Range("E30").Value = "  -0.93%  "

# Row 31
$ws.Range("D31").Value = "'0.0467"
$ws.Range("E31").Value = "  -2.52%  "

# Row 32
$ws.Range("D32").Value = "'3.29"
$ws.Range("E32").Value = "  -3.82%  "

# Row 33
$ws.Range("D33").Value = "1.348.41"
$ws.Range("E33").Value = "  -2.69%  "

# Row 34
$ws.Range("D34").Value = "'2.94"
$ws.Range("E34").Value = "  -4.02%  "

# Row 35
$ws.Range("E35").Value = "  -2.03%  "

# Row 36
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").Value = "'2.31"
$ws.Range("E36").Value = "  -0.67%  "

# Row 37
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").Value = "'0.956"
$ws.Range("E37").Value = "  -3.93%  "

# Row 38
$ws.Range("E38").Value = "  -2.45%  "

# Row 39
$ws.Range("D39").Value = "'0.537"
$ws.Range("E39").Value = "  -2.81%  "

# Row 40
$ws.Range("D40").Value = "'0.812"
$ws.Range("E40").Value = "  -3.49%  "

# Row 41
$ws.Range("E41").Value = "  +0.33%  "

# Row 42
$ws.Range("D42").Value = "'0.958"
$ws.Range("E42").Value = "  -3.43%  "

# Row 43
$ws.Range("D43").Value = "'5.33"
$ws.Range("E43").Value = "  -2.02%  "

# Row 44
$ws.Range("D44").Value = "'63.73"
$ws.Range("E44").Value = "  -2.53%  "

# Row 45
$ws.Range("D45").Value = "'1.74"
$ws.Range("E45").Value = "  -4.96%  "

# Row 46
$ws.Range("D46").Value = "1.730.14"
$ws.Range("E46").Value = "  -1.57%  "

# Row 47
$ws.Range("D47").Value = "'2.08"
$ws.Range("E47").Value = "  -3.26%  "

# Row 48
$ws.Range("D48").Value = "'87.22"
$ws.Range("E48").Value = "  -0.39%  "

# Row 49
$ws.Range("D49").Value = "0.0₇0994"
$ws.Range("E49").Value = "  -3.15%  "

# Row 50
$ws.Range("D50").Value = "'0.0965"
$ws.Range("E50").Value = "  -4.27%  "

# Row 51
$ws.Range("E51").Value = "  -1.08%  "
